$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at 190, pushing the existing rows 190-216 down to 194-220.
$ws.Rows("190:193").Insert()

# Populate the 4 newly inserted rows with the new weekly records.
# Columns A,B,C,E,F,G,H,I,J,R are constant boilerplate shared by every row in this block.
$rows = @(190, 191, 192, 193)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = 7
    $ws.Cells.Item($r, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
    $ws.Cells.Item($r, 3).Value = 'Ñuble'
    $ws.Cells.Item($r, 5).Value = 16
    $ws.Cells.Item($r, 6).Value = 'Fruta'
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = 'Cítricos'
    $ws.Cells.Item($r, 9).Value = 100102005
    $ws.Cells.Item($r, 10).Value = 'Naranja'
    $ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
}

# Row 190: Naranja / Cara cara / Primera
$ws.Range("D190").Value = 44449
$ws.Range("K190").Value = 'Cara cara'
$ws.Range("L190").Value = 'Primera'
$ws.Range("M190").Value = 240
$ws.Range("N190").Value = 5500
$ws.Range("O190").Value = 6000
$ws.Range("P190").Value = 5750
$ws.Range("Q190").Value = '$/malla 16 kilos'
$ws.Range("S190").Value = 359
$ws.Range("T190").Value = 16

# Row 191: Naranja / Cara cara / Segunda
$ws.Range("D191").Value = 44449
$ws.Range("K191").Value = 'Cara cara'
$ws.Range("L191").Value = 'Segunda'
$ws.Range("M191").Value = 180
$ws.Range("N191").Value = 4500
$ws.Range("O191").Value = 5000
$ws.Range("P191").Value = 4750
$ws.Range("Q191").Value = '$/malla 16 kilos'
$ws.Range("S191").Value = 297
$ws.Range("T191").Value = 16

# Row 192: Naranja / Navel Late / Primera
$ws.Range("D192").Value = 44449
$ws.Range("K192").Value = 'Navel Late'
$ws.Range("L192").Value = 'Primera'
$ws.Range("M192").Value = 240
$ws.Range("N192").Value = 5500
$ws.Range("O192").Value = 6000
$ws.Range("P192").Value = 5750
$ws.Range("Q192").Value = '$/bandeja 15 kilos granel'
$ws.Range("S192").Value = 383
$ws.Range("T192").Value = 15

# Row 193: Naranja / Navel Late / Segunda
$ws.Range("D193").Value = 44449
$ws.Range("K193").Value = 'Navel Late'
$ws.Range("L193").Value = 'Segunda'
$ws.Range("M193").Value = 190
$ws.Range("N193").Value = 4500
$ws.Range("O193").Value = 5000
$ws.Range("P193").Value = 4763
$ws.Range("Q193").Value = '$/bandeja 15 kilos granel'
$ws.Range("S193").Value = 318
$ws.Range("T193").Value = 15
